$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows("1:2").Insert()

$r1 = $ws.Range("A1:C1")
$r1.Value = "Unnamed: 0"
$ws.Range("A1").Value = "Unnamed: 0"
$ws.Range("B1").Value = "Unnamed: 1"
$ws.Range("C1").Value = "Unnamed: 2"

$r1.Font.Bold = $true
$r1.Borders.LineStyle = 1
$r1.HorizontalAlignment = -4108
$r1.VerticalAlignment = -4160

$ws.Range("A2").Value = "municipio"
$ws.Range("B2").Value = "CASOS"
$ws.Range("C2").Value = "ÓBITOS"

$ws.Range("A76").Value = "outros paises"
$ws.Range("B76").Value = 33

$ws.Range("A77").Value = "outros estados"
$ws.Range("B77").Value = 22
